$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), reusing the same formatting
# (bold, centered, thin box border) already applied to the other header
# cells (e.g. H1) by copying its format over before setting the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), rows 2-28
$data = @(
    @(1, 5),
    @(1, 5),
    @(4, 7),
    @(7, 9),
    @(3, 4),
    @(1, 3),
    @(8, 8),
    @(7, 9),
    @(4, 7),
    @(5, 6),
    @(8, 9),
    @(6, 9),
    @(8, 9),
    @(2, 3),
    @(6, 6),
    @(9, 9),
    @(5, 6),
    @(6, 7),
    @(6, 9),
    @(5, 6),
    @(5, 6),
    @(2, 6),
    @(7, 9),
    @(5, 9),
    @(6, 8),
    @(5, 7),
    @(4, 5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
